# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to Sheets/Aegis_Profits.xlsx per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 37372.93
$ws.Range("I64").Value = 69440.13
$ws.Range("J64").Value = 3015.2144
$ws.Range("K64").Value = 69440.13
$ws.Range("L64").Value = 3015.2144
$ws.Range("M64").Value = -69192.13
$ws.Range("N64").Value = -3511.2144
$ws.Range("H67").Value = 37372.93
$ws.Range("I67").Value = 69440.13
$ws.Range("J67").Value = 3015.2144
$ws.Range("K67").Value = 69440.13
$ws.Range("L67").Value = 3015.2144
$ws.Range("M67").Value = -68582.13
$ws.Range("N67").Value = -4731.2144
$ws.Range("H74").Value = 4997.5
$ws.Range("I74").Value = 4997.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4997.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4061.5
$ws.Range("N74").ClearContents()
$ws.Range("H76").Value = 3741.087
$ws.Range("I76").Value = 2999.2222
$ws.Range("J76").Value = 4218
$ws.Range("K76").Value = 2999.2222
$ws.Range("L76").Value = 4218
$ws.Range("M76").Value = -2684.2222
$ws.Range("N76").Value = -4848
$ws.Range("H77").Value = 4997.5
$ws.Range("I77").Value = 4997.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 24987.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -20307.5
$ws.Range("N77").ClearContents()
$ws.Range("H79").Value = 3741.087
$ws.Range("I79").Value = 2999.2222
$ws.Range("J79").Value = 4218
$ws.Range("K79").Value = 2999.2222
$ws.Range("L79").Value = 4218
$ws.Range("M79").Value = -1907.2222
$ws.Range("N79").Value = -6402
$ws.Range("H114").Value = 43722
$ws.Range("J114").Value = 43722
$ws.Range("L114").Value = 43722
$ws.Range("N114").Value = -52400
$ws.Range("H121").Value = 1950
$ws.Range("I121").Value = 2000
$ws.Range("K121").Value = 6000
$ws.Range("M121").Value = -4253
$ws.Range("H125").Value = 2765.08
$ws.Range("I125").Value = 2343.0588
$ws.Range("J125").Value = 3661.875
$ws.Range("K125").Value = 21087.5292
$ws.Range("L125").Value = 32956.875
$ws.Range("M125").Value = -18627.5292
$ws.Range("N125").Value = -37876.875
$ws.Range("H129").Value = 3258.1162
$ws.Range("J129").Value = 1013.4211
$ws.Range("L129").Value = 3040.2633
$ws.Range("N129").Value = -13040.2633
$ws.Range("H131").Value = 4116.7104
$ws.Range("I131").Value = 790.625
$ws.Range("J131").Value = 5003.6665
$ws.Range("K131").Value = 2371.875
$ws.Range("L131").Value = 15010.9995
$ws.Range("M131").Value = 2668.125
$ws.Range("N131").Value = -25090.9995
$ws.Range("H137").Value = 1678.909
$ws.Range("I137").Value = 1296.0769
$ws.Range("K137").Value = 3888.2307
$ws.Range("M137").Value = -1338.2307
$ws.Range("H141").Value = 1581.1428
$ws.Range("I141").Value = 1450.2
$ws.Range("J141").Value = 4200
$ws.Range("K141").Value = 4350.6
$ws.Range("L141").Value = 12600
$ws.Range("M141").Value = 829.3999999999996
$ws.Range("N141").Value = -22960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 8833.333000000001
$ws.Range("J15").Value = 8833.333000000001
$ws.Range("L15").Value = 8833.333000000001
$ws.Range("N15").Value = -9287.333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 163.33333
$ws.Range("I6").Value = 230
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 230
$ws.Range("L6").Value = 30
$ws.Range("M6").Value = -117
$ws.Range("N6").Value = -256
$ws.Range("H12").Value = 5000
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5340
$ws.Range("H58").Value = 6589.25
$ws.Range("I58").Value = 1396.5714
$ws.Range("J58").Value = 16502.545
$ws.Range("K58").Value = 1396.5714
$ws.Range("L58").Value = 16502.545
$ws.Range("M58").Value = -1193.5714
$ws.Range("N58").Value = -16908.545
$ws.Range("H134").Value = 1232.1818
$ws.Range("I134").Value = 1121.5555
$ws.Range("J134").Value = 1730
$ws.Range("K134").Value = 3364.6665
$ws.Range("L134").Value = 5190
$ws.Range("M134").Value = -829.6664999999998
$ws.Range("N134").Value = -10260
$ws.Range("H136").Value = 6589.25
$ws.Range("I136").Value = 1396.5714
$ws.Range("J136").Value = 16502.545
$ws.Range("K136").Value = 4189.7142
$ws.Range("L136").Value = 49507.63499999999
$ws.Range("M136").Value = -1639.7142
$ws.Range("N136").Value = -54607.63499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1316.6666
$ws.Range("I118").Value = 580
$ws.Range("J118").Value = 5000
$ws.Range("K118").Value = 1740
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = -497
$ws.Range("N118").Value = -17486
$ws.Range("H131").Value = 812.58
$ws.Range("I131").Value = 513.0769
$ws.Range("J131").Value = 857.3333
$ws.Range("K131").Value = 1539.2307
$ws.Range("L131").Value = 2571.9999
$ws.Range("M131").Value = 3500.7693
$ws.Range("N131").Value = -12651.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4545909
$ws.Range("J7").Value = 2502500
$ws.Range("L7").Value = 2502500
$ws.Range("N7").Value = -2502724
$ws.Range("H8").Value = 4545909
$ws.Range("J8").Value = 2502500
$ws.Range("L8").Value = 2502500
$ws.Range("N8").Value = -2502778
$ws.Range("H113").Value = 2533.5
$ws.Range("J113").Value = 1905.1111
$ws.Range("L113").Value = 1905.1111
$ws.Range("N113").Value = -6245.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2359.9333
$ws.Range("I7").Value = 1591
$ws.Range("M7").Value = -1479
$ws.Range("H16").Value = 46421.5
$ws.Range("I16").Value = 63019.938
$ws.Range("J16").Value = 2159
$ws.Range("K16").Value = 63019.938
$ws.Range("L16").Value = 2159
$ws.Range("M16").Value = -62849.938
$ws.Range("N16").Value = -2499
$ws.Range("H61").Value = 2594.1428
$ws.Range("J61").Value = 2985
$ws.Range("L61").Value = 2985
$ws.Range("N61").Value = -3389
$ws.Range("H113").Value = 2594.1428
$ws.Range("J113").Value = 2985
$ws.Range("L113").Value = 2985
$ws.Range("N113").Value = -7325
$ws.Range("H126").Value = 2359.9333
$ws.Range("I126").Value = 1591
$ws.Range("K126").Value = 4773
$ws.Range("M126").Value = -2303
$ws.Range("H127").Value = 41000
$ws.Range("J127").Value = 41000
$ws.Range("L127").Value = 41000
$ws.Range("N127").Value = -50920
$ws.Range("H136").Value = 1195.5641
$ws.Range("I136").Value = 1044.1765
$ws.Range("J136").Value = 2225
$ws.Range("K136").Value = 3132.5295
$ws.Range("L136").Value = 6675
$ws.Range("M136").Value = -582.5295000000001
$ws.Range("N136").Value = -11775

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 5007.625
$ws.Range("J49").Value = 5007.625
$ws.Range("L49").Value = 5007.625
$ws.Range("N49").Value = -5467.625
$ws.Range("H56").Value = 35803.8
$ws.Range("J56").Value = 38793.11
$ws.Range("L56").Value = 38793.11
$ws.Range("N56").Value = -40221.11
